$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 948.4286
$ws.Range("I2").Value = 1098.1666
$ws.Range("K2").Value = 1098.1666
$ws.Range("M2").Value = -985.1666
$ws.Range("H11").Value = 12.6875
$ws.Range("I11").Value = 12.6875
$ws.Range("K11").Value = 12.6875
$ws.Range("M11").Value = 127.3125
$ws.Range("H17").Value = 696.5484
$ws.Range("J17").Value = 713.13336
$ws.Range("L17").Value = 2139.40008
$ws.Range("N17").Value = -2475.40008
$ws.Range("H31").Value = 781
$ws.Range("I31").Value = 781
$ws.Range("K31").Value = 2343
$ws.Range("M31").Value = -2113
$ws.Range("H86").Value = 50003452
$ws.Range("I86").Value = 27780760
$ws.Range("J86").Value = 62503720
$ws.Range("K86").Value = 27780760
$ws.Range("L86").Value = 62503720
$ws.Range("M86").Value = -27779637
$ws.Range("N86").Value = -62505966
$ws.Range("H89").Value = 50003452
$ws.Range("I89").Value = 27780760
$ws.Range("J89").Value = 62503720
$ws.Range("K89").Value = 138903800
$ws.Range("L89").Value = 312518600
$ws.Range("M89").Value = -138898184
$ws.Range("N89").Value = -312529832
$ws.Range("H98").Value = 1023.6875
$ws.Range("I98").Value = 683.46155
$ws.Range("K98").Value = 683.46155
$ws.Range("M98").Value = 814.53845
$ws.Range("H122").Value = 1023.6875
$ws.Range("I122").Value = 683.46155
$ws.Range("K122").Value = 2050.38465
$ws.Range("M122").Value = 399.61535
$ws.Range("H132").Value = 1158.5
$ws.Range("I132").Value = 1245.326
$ws.Range("K132").Value = 3735.978
$ws.Range("M132").Value = -1205.978
$ws.Range("H138").Value = 1531.1177
$ws.Range("J138").Value = 2467.0908
$ws.Range("L138").Value = 7401.2724
$ws.Range("N138").Value = -17681.2724
$ws.Range("H141").Value = 9973.799999999999
$ws.Range("I141").Value = 9956.666999999999
$ws.Range("K141").Value = 29870.001
$ws.Range("M141").Value = -24690.001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 12502566
$ws.Range("J45").Value = 31252000
$ws.Range("L45").Value = 31252000
$ws.Range("N45").Value = -31252754
$ws.Range("H74").Value = 2285.3684
$ws.Range("I74").Value = 1899.0769
$ws.Range("K74").Value = 1899.0769
$ws.Range("M74").Value = -1025.0769
$ws.Range("H77").Value = 2285.3684
$ws.Range("I77").Value = 1899.0769
$ws.Range("K77").Value = 9495.3845
$ws.Range("M77").Value = -5127.3845
$ws.Range("H108").Value = 88996.664
$ws.Range("J108").Value = 88996.664
$ws.Range("L108").Value = 88996.664
$ws.Range("N108").Value = -96676.664
$ws.Range("H132").Value = 2028.8334
$ws.Range("I132").Value = 1591.7693
$ws.Range("K132").Value = 4775.3079
$ws.Range("M132").Value = -2245.3079

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 3000042.5
$ws.Range("I22").Value = 3713533.5
$ws.Range("K22").Value = 3713533.5
$ws.Range("M22").Value = -3713360.5
$ws.Range("H134").Value = 2000
$ws.Range("I134").Value = 2000
$ws.Range("K134").Value = 6000
$ws.Range("M134").Value = -3465

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 352.5484
$ws.Range("I7").Value = 238.46153
$ws.Range("J7").Value = 434.94446
$ws.Range("K7").Value = 238.46153
$ws.Range("L7").Value = 434.94446
$ws.Range("M7").Value = -125.46153
$ws.Range("N7").Value = -660.9444599999999
$ws.Range("H9").Value = 30569.715
$ws.Range("J9").Value = 30569.715
$ws.Range("L9").Value = 30569.715
$ws.Range("N9").Value = -30905.715
$ws.Range("H31").Value = 3113.182
$ws.Range("I31").Value = 1899.375
$ws.Range("J31").Value = 6350
$ws.Range("K31").Value = 1899.375
$ws.Range("L31").Value = 6350
$ws.Range("M31").Value = -1604.375
$ws.Range("N31").Value = -6940
$ws.Range("H34").Value = 3113.182
$ws.Range("I34").Value = 1899.375
$ws.Range("J34").Value = 6350
$ws.Range("K34").Value = 1899.375
$ws.Range("L34").Value = 6350
$ws.Range("M34").Value = -1697.375
$ws.Range("N34").Value = -6754
$ws.Range("H122").Value = 2055.6956
$ws.Range("I122").Value = 1640.0588
$ws.Range("K122").Value = 4920.1764
$ws.Range("M122").Value = -2470.1764
$ws.Range("H132").Value = 1894.7826
$ws.Range("I132").Value = 1674.7142
$ws.Range("J132").Value = 2237.111
$ws.Range("K132").Value = 5024.142599999999
$ws.Range("L132").Value = 6711.333
$ws.Range("M132").Value = -2494.142599999999
$ws.Range("N132").Value = -11771.333

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1004.3333
$ws.Range("I5").Value = 849.9
$ws.Range("J5").Value = 1197.375
$ws.Range("K5").Value = 2549.7
$ws.Range("L5").Value = 3592.125
$ws.Range("M5").Value = -2437.7
$ws.Range("N5").Value = -3816.125
$ws.Range("H135").Value = 1004.3333
$ws.Range("I135").Value = 849.9
$ws.Range("J135").Value = 1197.375
$ws.Range("K135").Value = 7649.099999999999
$ws.Range("L135").Value = 10776.375
$ws.Range("M135").Value = -5114.099999999999
$ws.Range("N135").Value = -15846.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 250002500
$ws.Range("I80").Value = 500001540
$ws.Range("J80").Value = 3442
$ws.Range("K80").Value = 500001540
$ws.Range("L80").Value = 3442
$ws.Range("M80").Value = -500000542
$ws.Range("N80").Value = -5438
$ws.Range("H83").Value = 250002500
$ws.Range("I83").Value = 500001540
$ws.Range("J83").Value = 3442
$ws.Range("K83").Value = 2500007700
$ws.Range("L83").Value = 17210
$ws.Range("M83").Value = -2500002708
$ws.Range("N83").Value = -27194
$ws.Range("H102").Value = 2290.2222
$ws.Range("I102").Value = 2290.2222
$ws.Range("K102").Value = 2290.2222
$ws.Range("M102").Value = -668.2222000000002
$ws.Range("H132").Value = 3250.5527
$ws.Range("I132").Value = 2353.1482
$ws.Range("J132").Value = 5453.273
$ws.Range("K132").Value = 7059.444600000001
$ws.Range("L132").Value = 16359.819
$ws.Range("M132").Value = -4529.444600000001
$ws.Range("N132").Value = -21419.819

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 2713.75
$ws.Range("J55").Value = 8250
$ws.Range("L55").Value = 8250
$ws.Range("N55").Value = -8596
$ws.Range("H62").Value = 13000
$ws.Range("J62").Value = 15000
$ws.Range("L62").Value = 15000
$ws.Range("N62").Value = -16248
$ws.Range("H65").Value = 13000
$ws.Range("J65").Value = 15000
$ws.Range("L65").Value = 45000
$ws.Range("N65").Value = -51240
$ws.Range("H76").Value = 12070.75
$ws.Range("J76").Value = 16250
$ws.Range("L76").Value = 16250
$ws.Range("N76").Value = -16926
$ws.Range("H79").Value = 12070.75
$ws.Range("J79").Value = 16250
$ws.Range("L79").Value = 16250
$ws.Range("N79").Value = -18590
$ws.Range("H136").Value = 1704.28
$ws.Range("I136").Value = 1508.7273
$ws.Range("J136").Value = 1857.9286
$ws.Range("K136").Value = 4526.1819
$ws.Range("L136").Value = 5573.7858
$ws.Range("M136").Value = -1976.1819
$ws.Range("N136").Value = -10673.7858

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 10580.5
$ws.Range("J41").Value = 10652.4
$ws.Range("L41").Value = 10652.4
$ws.Range("N41").Value = -11432.4
$ws.Range("H63").Value = 9489.799999999999
$ws.Range("I63").Value = 1200
$ws.Range("J63").Value = 11562.25
$ws.Range("K63").Value = 1200
$ws.Range("L63").Value = 11562.25
$ws.Range("M63").Value = -576
$ws.Range("N63").Value = -12810.25
$ws.Range("H66").Value = 9489.799999999999
$ws.Range("I66").Value = 1200
$ws.Range("J66").Value = 11562.25
$ws.Range("K66").Value = 3600
$ws.Range("L66").Value = 34686.75
$ws.Range("M66").Value = -480
$ws.Range("N66").Value = -40926.75
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
$ws.Range("H86").Value = 55000
$ws.Range("J86").Value = 55000
$ws.Range("L86").Value = 55000
$ws.Range("N86").Value = -57246
$ws.Range("H89").Value = 55000
$ws.Range("J89").Value = 55000
$ws.Range("L89").Value = 275000
$ws.Range("N89").Value = -286232
$ws.Range("H107").Value = 1650.9166
$ws.Range("I107").Value = 699.8333
$ws.Range("J107").Value = 2602
$ws.Range("K107").Value = 2099.4999
$ws.Range("L107").Value = 7806
$ws.Range("M107").Value = -179.4998999999998
$ws.Range("N107").Value = -11646
$ws.Range("H132").Value = 1505.3667
$ws.Range("I132").Value = 1138.3334
$ws.Range("K132").Value = 3415.0002
$ws.Range("M132").Value = -885.0001999999999
$ws.Range("H136").Value = 1864.7142
$ws.Range("I136").Value = 1864.7142
$ws.Range("K136").Value = 5594.142599999999
$ws.Range("M136").Value = -3044.142599999999
